$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.800385475158691
$ws.Range("B1").Value = -1
$ws.Range("C1").Value = 2.674353837966919
$ws.Range("D1").Value = 1.128020048141479
$ws.Range("E1").Value = 0.7688800692558289
